$wb = $excel.ActiveWorkbook

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $last)
$ws.Name = "ODI Batting Extra"

# Header row
$ws.Cells.Item(1,1).Value = "MATCH_CODE"
$ws.Cells.Item(1,2).Value = "BATTING_POSITION"
$ws.Cells.Item(1,3).Value = "NUM_4"
$ws.Cells.Item(1,4).Value = "NUM_6"
$ws.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Cells.Item(1,6).Value = "MAN_OF_MATCH"
$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("A1:F1").Borders.LineStyle = 1
$ws.Range("A1:F1").HorizontalAlignment = -4108
$ws.Range("A1:F1").VerticalAlignment = -4160

# Data rows. Columns A, C, D, E hold text that looks numeric (match codes,
# counts, percentages) in the source data, so each such cell is switched to
# text format ("@") right before its value is written, to avoid Excel
# auto-converting the literal to a real number/percent. Truly-empty cells
# in these rows are left untouched (no cell is created), matching the
# source data which leaves those fields blank.
$ws.Range("A2").NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "4634"
$ws.Cells.Item(2,6).Value = "NO"

$ws.Range("A3").NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "4638"
$ws.Cells.Item(3,2).Value = 5
$ws.Range("C3:E3").NumberFormat = "@"
$ws.Cells.Item(3,3).Value = "5"
$ws.Cells.Item(3,4).Value = "2"
$ws.Cells.Item(3,5).Value = "26.18%"
$ws.Cells.Item(3,6).Value = "NO"

$ws.Range("A4").NumberFormat = "@"
$ws.Cells.Item(4,1).Value = "4641"
$ws.Cells.Item(4,2).Value = 4
$ws.Range("C4:E4").NumberFormat = "@"
$ws.Cells.Item(4,3).Value = "1"
$ws.Cells.Item(4,4).Value = "0"
$ws.Cells.Item(4,5).Value = "11.65%"
$ws.Cells.Item(4,6).Value = "NO"

$ws.Range("A5").NumberFormat = "@"
$ws.Cells.Item(5,1).Value = "4686"
$ws.Cells.Item(5,6).Value = "NO"

$ws.Range("A6").NumberFormat = "@"
$ws.Cells.Item(6,1).Value = "4688"
$ws.Cells.Item(6,2).Value = 6
$ws.Range("C6:E6").NumberFormat = "@"
$ws.Cells.Item(6,3).Value = "2"
$ws.Cells.Item(6,4).Value = "0"
$ws.Cells.Item(6,5).Value = "13.74%"
$ws.Cells.Item(6,6).Value = "NO"

$ws.Range("A7").NumberFormat = "@"
$ws.Cells.Item(7,1).Value = "4690"
$ws.Cells.Item(7,6).Value = "NO"

Write-Host "done"
